# Robetroid Course of Action.docx
# Commit: "Added Player Hit animation with Knockback"
#
# The underlying XML diff for this edit only touches highlight colors
# (yellow/green -> cyan) on a cluster of bullet paragraphs describing the
# Player's "Hurt" animation and the "Enemies (Robots)" bullet, plus the
# auto-managed "_GoBack" bookmark (Word's "last edit location" marker)
# shifting from the end of the "Hurt" bullet down to the end of the
# "Hurts Robotroid when touching him" bullet - consistent with the author
# having made their last edit there while adding the new Knockback note.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the unique paragraph run:
#      "Add Animations" / "Hurt" / "HUD" / "# of lives" / "Health bar" /
#      "Enemies (Robots)" / "Moves along the platform it's on" /
#      "Hurts Robotroid when touching him"
#    (the document has more than one "Add Animations" bullet, so match on
#    the whole sequence to get the right one).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$anchor = -1
for ($i = 1; $i -le $count - 7; $i++) {
    $t1 = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    $t2 = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd([char]13)
    $t3 = $d.Paragraphs.Item($i + 2).Range.Text.TrimEnd([char]13)
    if ($t1 -eq "Add Animations" -and $t2 -eq "Hurt" -and $t3 -eq "HUD") {
        $anchor = $i
        break
    }
}

if ($anchor -eq -1) {
    throw "Could not locate the target 'Add Animations' / 'Hurt' bullet cluster"
}

$pAddAnimations = $d.Paragraphs.Item($anchor)
$pHurt          = $d.Paragraphs.Item($anchor + 1)
$pHud           = $d.Paragraphs.Item($anchor + 2)
$pLives         = $d.Paragraphs.Item($anchor + 3)
$pHealthBar     = $d.Paragraphs.Item($anchor + 4)
$pEnemies       = $d.Paragraphs.Item($anchor + 5)
$pMoves         = $d.Paragraphs.Item($anchor + 6)
$pHurtsRobot    = $d.Paragraphs.Item($anchor + 7)

# ------------------------------------------------------------------
# 2. Re-highlight "Add Animations" and "Hurt" from yellow to cyan
#    (wdTurquoise = 3). Using Range.Font.HighlightColorIndex (rather than
#    Range.HighlightColorIndex) so both the run text *and* the paragraph
#    mark's run properties pick up the new colour, matching the diff.
# ------------------------------------------------------------------
$pAddAnimations.Range.Font.HighlightColorIndex = 3
$pHurt.Range.Font.HighlightColorIndex = 3

# ------------------------------------------------------------------
# 3. Re-highlight "Enemies (Robots)" from green to cyan.
# ------------------------------------------------------------------
$pEnemies.Range.Font.HighlightColorIndex = 3

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark (Word's invisible "last edit" marker)
#    from right after "Hurt" to right after "...when touching him",
#    i.e. reflecting that the author's final touch for this commit was
#    on that bullet. Adding a bookmark named "_GoBack" again relocates
#    the existing one rather than creating a duplicate.
# ------------------------------------------------------------------
$goBackTarget = $pHurtsRobot.Range.End - 1
$goBackRange = $d.Range($goBackTarget, $goBackTarget)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "Updated highlights on paragraphs $anchor..$($anchor+7) and relocated _GoBack bookmark."
